$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "NAQUICHE SILVA MARIA LIZETH"
$ws.Range("B2").Value = 106

$ws.Range("A3").Value = "NAQUICHE MECHATO ENMA"
$ws.Range("B3").Value = 101

$ws.Range("A4").Value = "CUBAS GARCIA ROSA ANITA"
$ws.Range("B4").Value = 76

$ws.Range("A5").Value = "MANOSALVA RUIZ SANDRA KAROLINE"
$ws.Range("B5").Value = 73

$ws.Range("A6").Value = "CORAS QUISPE JORGE AMERICO"
$ws.Range("B6").Value = 72

$ws.Range("A7").Value = "BECERRA ASMAT CAROL STEFANY"
$ws.Range("B7").Value = 61

$ws.Range("A8").Value = "PACHECO ALISON"
$ws.Range("B8").Value = 50

$ws.Range("A9").Value = "SAUCEDO CABRERA CARLOS ALEXANDER"
$ws.Range("B9").Value = 46

$ws.Range("A10").Value = "CASTREJON TELLO GRECIA"
$ws.Range("B10").Value = 20
